# Update the "addListItem" sheet: rename the list item from "LinuxAK" to "ListLin"
$wb = $excel.ActiveWorkbook

$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "ListLin"

# Update the "createUser" sheet: bump the test user id from 2721 to 2724
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 2724
